$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both `"Barbie`" and `"Oppenheimer`" has been recorded successfully.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been finalized.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("D4").Value = "Barbie_was_selected, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Oppenheimer`" for the movie to be shown on Friday.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to show on Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been finalized with no movie selected for Friday.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Oppenheimer`" will be acquired for the showing on Friday.`n"
$ws.Range("D8").Value = "Oppenheimer_was_selected, "
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has resulted in no selection.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday was not reached, leading to no consensus.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie has been selected for Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been made to not acquire any movie for the Friday show, as there was no agreement reached during the discussion.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: I have recorded the decision as `"no_decision`" since the committee did not reach a clear agreement about which movie to show on Friday.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday could not be made, as there was no consensus reached by the committee.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: I have recorded the decision as `"no decision`" about the movie to be shown on Friday.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded and the rights for `"Barbie`" have been successfully acquired.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision being made about the movie for Friday.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday remains unresolved.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no consensus.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision has been successfully recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("D24").Value = "Barbie_was_selected, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights to `"Barbie.`"`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The function for no decision has been successfully called, indicating that the committee did not reach a conclusion regarding the movie selection for Friday.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be shown on Friday.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday was not reached, leading to a no_decision outcome.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: No movie was selected in this meeting.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded: no decision was made regarding the movie to show on Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made explicitly, so I have recorded that there is no decision at this time.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded, and there was no clear outcome regarding which movie to show on Friday.`n"
$ws.Range("D35").Value = "no_decision, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: I have recorded the decision as no decision was made regarding which movie will be shown on Friday.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision process for Friday's movie did not result in a clear choice, indicating no agreement was reached. If further discussion is required in the future, that would be an appropriate next step.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded.`n"
$ws.Range("D39").Value = "both_movies, "
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision process concluded without a consensus, resulting in no movie being selected for Friday.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday was not made during the committee discussion.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The function has been successfully called, indicating that there was no decision made about which movie to show.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no consensus was reached regarding the movie to be shown on Friday.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie for Friday.`n"
$ws.Range("D45").Value = "no_decision, "
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no agreement.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision about the movie for Friday has been recorded as no decision.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie has been selected for Friday.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded: `"Oppenheimer`" will be shown on Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday was not finalized.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("D53").Value = "both_movies, "
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision process has concluded without a definitive choice for Friday's movie, leading to no decision being made.`n"
$ws.Range("D54").Value = "no_decision, "
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been recorded as no consensus was reached regarding the movie selection for Friday.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision to select `"Barbie`" as the movie for Friday has been recorded successfully.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie has been selected for Friday.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" based on the agreement reached in the discussion.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected in the meeting.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday has resulted in no conclusion.`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no choice of a movie for Friday's screening was made during the discussion.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made, so I confirm that there is no decision.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision-making process concluded without an agreement on which movie to show on Friday, resulting in no decision being made.`n"
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to both movies.`n"
$ws.Range("C68").Value = "MSG: None`n`nMSG: The decision process concluded without a definitive choice for Friday's movie, resulting in no decision being made.`n"

# Restore automatic row heights so the row-height metadata matches the
# original file (Excel auto-expands row height when multi-line text is
# written into a cell via COM automation).
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).AutoFit()
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(20).AutoFit()
$ws.Rows.Item(21).AutoFit()
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(25).AutoFit()
$ws.Rows.Item(26).AutoFit()
$ws.Rows.Item(27).AutoFit()
$ws.Rows.Item(28).AutoFit()
$ws.Rows.Item(29).AutoFit()
$ws.Rows.Item(30).AutoFit()
$ws.Rows.Item(31).AutoFit()
$ws.Rows.Item(32).AutoFit()
$ws.Rows.Item(33).AutoFit()
$ws.Rows.Item(34).AutoFit()
$ws.Rows.Item(35).AutoFit()
$ws.Rows.Item(36).AutoFit()
$ws.Rows.Item(37).AutoFit()
$ws.Rows.Item(38).AutoFit()
$ws.Rows.Item(39).AutoFit()
$ws.Rows.Item(40).AutoFit()
$ws.Rows.Item(41).AutoFit()
$ws.Rows.Item(42).AutoFit()
$ws.Rows.Item(43).AutoFit()
$ws.Rows.Item(44).AutoFit()
$ws.Rows.Item(45).AutoFit()
$ws.Rows.Item(46).AutoFit()
$ws.Rows.Item(47).AutoFit()
$ws.Rows.Item(48).AutoFit()
$ws.Rows.Item(49).AutoFit()
$ws.Rows.Item(50).AutoFit()
$ws.Rows.Item(51).AutoFit()
$ws.Rows.Item(52).AutoFit()
$ws.Rows.Item(53).AutoFit()
$ws.Rows.Item(54).AutoFit()
$ws.Rows.Item(55).AutoFit()
$ws.Rows.Item(56).AutoFit()
$ws.Rows.Item(57).AutoFit()
$ws.Rows.Item(58).AutoFit()
$ws.Rows.Item(59).AutoFit()
$ws.Rows.Item(60).AutoFit()
$ws.Rows.Item(61).AutoFit()
$ws.Rows.Item(62).AutoFit()
$ws.Rows.Item(63).AutoFit()
$ws.Rows.Item(64).AutoFit()
$ws.Rows.Item(65).AutoFit()
$ws.Rows.Item(66).AutoFit()
$ws.Rows.Item(67).AutoFit()
$ws.Rows.Item(68).AutoFit()
